$d = $word.ActiveDocument

# Replace "+0" with "0" everywhere in the document body
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result = $find.Execute("+0", $true, $false, $false, $false, $false, $true, 1, $false, "0", 2)
Write-Output "Replace +0 -> 0 result: $result"

# Replace "MONTH: TOTAL" with "MONTH: JUNE-OCT"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$result2 = $find2.Execute("MONTH: TOTAL", $true, $false, $false, $false, $false, $true, 1, $false, "MONTH: JUNE-OCT", 2)
Write-Output "Replace MONTH result: $result2"
